# This script adds a new column BB (a copy of column BA, the most recent
# QoQ forecast column) reflecting an updated forecast run that now also
# folds in an EQUIPMENT evaluation, and appends one additional forecast
# row (row 83) for the next quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate column BA (header date + all forecast rows) into the new
#    column BB, carrying over values *and* styles (date formatting on
#    row 1, borders, etc.) in one shot.
$ws.Range("BA1:BA82").Copy($ws.Range("BB1:BB82"))

# 2) The new column's header date is one quarter later than BA's.
$ws.Range("BB1").Value = 45986

# 3) Starting at row 71 (2025-06-30) the forecast was revised because the
#    EQUIPMENT evaluation is now included, so overwrite those cells with
#    the newly calculated values.
$ws.Range("BB71").Value = -0.1118837721692358
$ws.Range("BB72").Value = 0.3266766184601977
$ws.Range("BB73").Value = 0.325608361860148
$ws.Range("BB74").Value = 0.325608361860148
$ws.Range("BB75").Value = 0.325608361860148
$ws.Range("BB76").Value = 0.325608361860148
$ws.Range("BB77").Value = 0.325608361860148
$ws.Range("BB78").Value = 0.325608361860148
$ws.Range("BB79").Value = 0.325608361860148
$ws.Range("BB80").Value = 0.325608361860148
$ws.Range("BB81").Value = 0.325608361860148
$ws.Range("BB82").Value = 0.325608361860148

# 4) Append one more quarter of data as row 83: the date in column A
#    (copied so it keeps the same date-number style as the rest of
#    column A) and the forecast value in the new column BB.
$ws.Range("A82").Copy($ws.Range("A83"))
$ws.Range("A83").Value = 46934

$ws.Range("BA82").Copy($ws.Range("BB83"))
$ws.Range("BB83").Value = 0.325608361860148
